$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old "Mesa" column (D) is no longer used; drop it, then insert a new
# column before "Colegio" (currently C) to hold the new "Email" data. This
# shifts "Colegio" from C to D, matching the target layout.
$ws.Columns.Item(4).Delete()
$ws.Columns.Item(3).Insert()

# New "Email" column header + values.
$ws.Range("C1").Value = "Email"
$ws.Range("C2").Value = "a@gmail.com"
$ws.Range("C3").Value = "b@gmail.com"
$ws.Range("C4").Value = "c@gmail.com"

# Turn the email cells into mailto: hyperlinks (this also applies the
# built-in "Hyperlink" style/font to each cell).
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:a@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:b@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:c@gmail.com")

# Move the active selection to D6, as in the edited workbook.
[void]$ws.Range("D6").Select()
